$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates ---
$ws.Range("A8").Value = "Volume 30   Number  18"
$ws.Range("C9").Value = "Report Covering the Week  5/1/2023  Through  5/7/2023"

# --- Crime statistics table updates ---
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0'
$ws.Range("A14").Copy()
$ws.Range("D14").PasteSpecial(-4122)
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '***.*'
$ws.Range("A14").Copy()
$ws.Range("E14").PasteSpecial(-4122)
$ws.Range("A15").Value = 'Rape'
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = '0'
$ws.Range("A14").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("E15").Value = '***.*'
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = 0
$ws.Range("N15").Value = -25
$ws.Range("F16").Value = 13
$ws.Range("G16").Value = 6
$ws.Range("H16").Value = 116.666666666667
$ws.Range("I16").Value = 45
$ws.Range("J16").Value = 25
$ws.Range("K16").Value = 80
$ws.Range("L16").Value = 55.172413793103
$ws.Range("M16").Value = -2.173913043478
$ws.Range("N16").Value = -76.804123711340
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 17
$ws.Range("G17").Value = 18
$ws.Range("H17").Value = -5.555555555555
$ws.Range("I17").Value = 74
$ws.Range("J17").Value = 83
$ws.Range("K17").Value = -10.843373493975
$ws.Range("L17").Value = -27.450980392156
$ws.Range("M17").Value = 27.586206896551
$ws.Range("N17").Value = -44.776119402985
$ws.Range("G14").Copy()
$ws.Range("C18").PasteSpecial(-4122)
$ws.Range("C18").Value = 2
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0'
$ws.Range("A14").Copy()
$ws.Range("D18").PasteSpecial(-4122)
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '***.*'
$ws.Range("A14").Copy()
$ws.Range("E18").PasteSpecial(-4122)
$ws.Range("F18").Value = 5
$ws.Range("G18").Value = 4
$ws.Range("H18").Value = 25
$ws.Range("I18").Value = 26
$ws.Range("K18").Value = 4
$ws.Range("L18").Value = 62.5
$ws.Range("M18").Value = -29.729729729729
$ws.Range("N18").Value = -89.166666666666
$ws.Range("C19").Value = 2
$ws.Range("D19").Value = 4
$ws.Range("E19").Value = -50
$ws.Range("F19").Value = 11
$ws.Range("G19").Value = 17
$ws.Range("H19").Value = -35.294117647058
$ws.Range("I19").Value = 57
$ws.Range("J19").Value = 70
$ws.Range("K19").Value = -18.571428571428
$ws.Range("L19").Value = 14
$ws.Range("M19").Value = 72.727272727272
$ws.Range("N19").Value = -43
$ws.Range("F20").Value = 3
$ws.Range("H20").Value = -25
$ws.Range("I20").Value = 23
$ws.Range("J20").Value = 13
$ws.Range("K20").Value = 76.923076923076
$ws.Range("L20").Value = 21.052631578947
$ws.Range("M20").Value = -14.814814814814
$ws.Range("N20").Value = -87.894736842105
$ws.Range("D21").Value = 12
$ws.Range("E21").Value = 8.333333333333
$ws.Range("F21").Value = 50
$ws.Range("G21").Value = 51
$ws.Range("H21").Value = -1.960784313725
$ws.Range("I21").Value = 231
$ws.Range("J21").Value = 227
$ws.Range("K21").Value = 1.762114537444
$ws.Range("L21").Value = 2.212389380530
$ws.Range("M21").Value = 9.478672985781
$ws.Range("N21").Value = -73.478760045924
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0'
$ws.Range("A14").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '***.*'
$ws.Range("A14").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("F22").NumberFormat = "@"
$ws.Range("F22").Value = '0'
$ws.Range("A14").Copy()
$ws.Range("F22").PasteSpecial(-4122)
$ws.Range("G22").Value = 1
$ws.Range("H22").Value = -100
$ws.Range("N22").Value = '***.*'
$ws.Range("G14").Copy()
$ws.Range("C23").PasteSpecial(-4122)
$ws.Range("C23").Value = 1
$ws.Range("G14").Copy()
$ws.Range("D23").PasteSpecial(-4122)
$ws.Range("D23").Value = 2
$ws.Range("H14").Copy()
$ws.Range("E23").PasteSpecial(-4122)
$ws.Range("E23").Value = -50
$ws.Range("F23").Value = 4
$ws.Range("G23").Value = 2
$ws.Range("H23").Value = 100
$ws.Range("I23").Value = 30
$ws.Range("J23").Value = 23
$ws.Range("K23").Value = 30.434782608695
$ws.Range("L23").Value = -9.090909090909
$ws.Range("M23").Value = 87.5
$ws.Range("N23").Value = '***.*'
$ws.Range("C24").Value = 20
$ws.Range("D24").Value = 14
$ws.Range("E24").Value = 42.857142857142
$ws.Range("F24").Value = 60
$ws.Range("G24").Value = 46
$ws.Range("H24").Value = 30.434782608695
$ws.Range("I24").Value = 194
$ws.Range("J24").Value = 175
$ws.Range("K24").Value = 10.857142857142
$ws.Range("L24").Value = 3.743315508021
$ws.Range("M24").Value = 61.666666666666
$ws.Range("N24").Value = '***.*'
$ws.Range("C25").Value = 6
$ws.Range("D25").Value = 11
$ws.Range("E25").Value = -45.454545454545
$ws.Range("F25").Value = 24
$ws.Range("G25").Value = 33
$ws.Range("H25").Value = -27.272727272727
$ws.Range("I25").Value = 125
$ws.Range("J25").Value = 123
$ws.Range("K25").Value = 1.626016260162
$ws.Range("L25").Value = -17.763157894736
$ws.Range("M25").Value = -8.088235294117
$ws.Range("N25").Value = '***.*'
$ws.Range("E26").Value = '***.*'
$ws.Range("F26").Value = 2
$ws.Range("G26").Value = 1
$ws.Range("H26").Value = 100
$ws.Range("I26").Value = 13
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 8.333333333333
$ws.Range("M26").Value = '***.*'
$ws.Range("N26").Value = '***.*'
$ws.Range("G14").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("C27").Value = 2
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0'
$ws.Range("A14").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '***.*'
$ws.Range("A14").Copy()
$ws.Range("E27").PasteSpecial(-4122)
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 14
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = -6.666666666666
$ws.Range("M27").Value = '***.*'
$ws.Range("N27").Value = '***.*'
$ws.Range("E28").Value = '***.*'
$ws.Range("G28").Value = 1
$ws.Range("H28").Value = 0
$ws.Range("M28").Value = -33.333333333333
$ws.Range("N28").Value = -62.5
$ws.Range("E29").Value = '***.*'
$ws.Range("G29").Value = 1
$ws.Range("H29").Value = 0
$ws.Range("M29").Value = -55.555555555555
$ws.Range("N29").Value = -75
$ws.Range("E30").Value = '***.*'
$ws.Range("F30").Value = 1
$ws.Range("H30").Value = 0
$ws.Range("M30").Value = '***.*'
$ws.Range("N30").Value = '***.*'
$ws.Range("A37").Value = 'Rape'

$excel.CutCopyMode = 0

